# S06/G02: Risk management engine (limits and overrides)
#
# Applies the content/status updates described by the commit, matching the
# target XML diff:
#  - Row 18 (S02_G03_TF003): rewrite deviations/remarks/pending-work text to
#    reflect that the Settings page grew beyond read-only.
#  - Rows 43-45 (S06_G01 tasks): pick up the same "normal" direct formatting
#    already used by the surrounding data rows (no wrap, bottom-aligned)
#    for columns F/H/I, without touching their text.
#  - Rows 46-49 (S06_G02 tasks): fill in the previously-empty
#    deviations/remarks/pending-work columns and flip status from
#    "pending" to "implemented", using the same normal formatting as the
#    rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell's text and make sure it uses the same plain
# (non-wrapped, bottom-aligned) direct formatting used elsewhere in the
# F/H/I "notes" columns of this sheet.
function Set-NoteCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $cell.WrapText = $false
    $cell.VerticalAlignment = -4107
    $cell.HorizontalAlignment = 1
}

# Helper: normalize an existing cell's direct formatting only (value is
# left untouched).
function Set-NoteFormat($addr) {
    $cell = $ws.Range($addr)
    $cell.WrapText = $false
    $cell.VerticalAlignment = -4107
    $cell.HorizontalAlignment = 1
}

# --- Row 18: S02_G03_TF003 (text only; formatting (s="5") is unchanged) ---
$ws.Range("F18").Value = "Settings page started as read-only and was later extended to allow editing strategy execution modes and creating GLOBAL or per-strategy risk settings directly from the UI."
$ws.Range("H18").Value = "Strategies and risk_settings can now be inspected and, for key fields, edited from a single Settings screen (mode toggle and risk creation)."
$ws.Range("I18").Value = "Add delete/edit flows for existing risk rows and finer-grained admin controls once configuration needs grow."

# --- Row 43: S06_G01_TB001 (text unchanged, formatting normalized) --------
Set-NoteFormat "F43"
Set-NoteFormat "H43"
Set-NoteFormat "I43"

# --- Row 44: S06_G01_TB002 (text unchanged, formatting normalized) --------
Set-NoteFormat "F44"
Set-NoteFormat "H44"
Set-NoteFormat "I44"

# --- Row 45: S06_G01_TF003 (text unchanged, formatting normalized) --------
Set-NoteFormat "F45"
Set-NoteFormat "H45"
Set-NoteFormat "I45"

# --- Row 46: S06_G02_TB001 --------------------------------------------------
Set-NoteCell "F46" "Implemented a risk evaluation service that uses RiskSettings (GLOBAL and per-strategy) to enforce max_quantity_per_order, max_order_value, and allow_short_selling before broker calls."
$ws.Range("G46").Value = "implemented"
Set-NoteCell "H46" "Risk checks run inside the order execution path for both MANUAL queue executes and AUTO strategy orders."
Set-NoteCell "I46" "Extend the risk engine to incorporate max_daily_loss and max_open_positions once realized PnL and positions are tracked in later sprints."

# --- Row 47: S06_G02_TB002 --------------------------------------------------
Set-NoteCell "F47" "Orders violating hard limits are now blocked with status REJECTED_RISK and a human-readable explanation; when clamp_mode=CLAMP, quantities are reduced instead of rejected where possible."
$ws.Range("G47").Value = "implemented"
Set-NoteCell "H47" "Risk decisions are stored on the Order via status/error_message and are visible in the Orders history UI."
Set-NoteCell "I47" "Refine risk messages and add per-rule identification if we need more granular auditing later."

# --- Row 48: S06_G02_TB003 --------------------------------------------------
Set-NoteCell "F48" "Risk checks are invoked inside the shared order execution endpoint, so both AUTO strategy executions (from the webhook) and MANUAL queue executes are evaluated before any broker call."
$ws.Range("G48").Value = "implemented"
Set-NoteCell "H48" "AUTO and MANUAL flows now share the same risk gate in execute_order; broker requests are only sent if the risk engine allows the order."
Set-NoteCell "I48" "As we add more execution paths (e.g., bulk actions), ensure they all call through the same risk-aware execution helper."

# --- Row 49: S06_G02_TF004 --------------------------------------------------
Set-NoteCell "F49" "Risk-related rejections and clamps are surfaced via the existing Error column in the Orders UI, showing detailed messages from the risk engine."
$ws.Range("G49").Value = "implemented"
Set-NoteCell "H49" "Users can see which orders were blocked (REJECTED_RISK) or had their quantities adjusted before execution."
Set-NoteCell "I49" "Consider adding explicit risk badges or tooltips in a later UX-focused sprint to differentiate risk notes from broker errors."
